$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("E2").ClearContents()
$ws.Range("E2").Select()
